$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 405, shifting existing rows 405-438 down to 406-439.
$ws.Rows.Item(405).Insert()

# Populate the newly inserted row 405 with its data.
$ws.Cells.Item(405, 1).Value = 3
$ws.Cells.Item(405, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(405, 3).Value = "Coquimbo"
$ws.Cells.Item(405, 4).Value = 44783
$ws.Cells.Item(405, 5).Value = 5
$ws.Cells.Item(405, 6).Value = 100112017
$ws.Cells.Item(405, 7).Value = "Apio"
$ws.Cells.Item(405, 8).Value = "Americana (o)"
$ws.Cells.Item(405, 9).Value = "Primera"
$ws.Cells.Item(405, 10).Value = 230
$ws.Cells.Item(405, 11).Value = 9000
$ws.Cells.Item(405, 12).Value = 9500
$ws.Cells.Item(405, 13).Value = 9261
$ws.Cells.Item(405, 14).Value = "$/docena de matas"
$ws.Cells.Item(405, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(405, 16).Value = 1544
$ws.Cells.Item(405, 17).Value = 6
$ws.Cells.Item(405, 18).Value = "Hortaliza"
